# Edit: "group power plants per technology and installation year"
#
# Changes applied:
#  - Coupling Parameters!B18: investment_initialization_years 3 -> 0
#  - Coupling Parameters!B37: minimal_last_years_IRR "NOTSET" -> TRUE (boolean)
#  - Coupling Parameters: new row 40 -> A40 "groups power plants per installed year", B40 = TRUE
#  - Coupling Parameters sheet view: scroll to A18, select B38

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Coupling Parameters")

# investment_initialization_years: 3 -> 0
$ws.Range("B18").Value = 0

# minimal_last_years_IRR: "NOTSET" -> TRUE
$ws.Range("B37").Value = $true

# New row: groups power plants per installed year -> TRUE
$ws.Range("A40").Value = "groups power plants per installed year"
$ws.Range("B40").Value = $true

# Update view: scroll position and selection on the Coupling Parameters sheet
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 18
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B38").Select()

$wb.Save()
